{"js": "// Remove the trailing \"Ver no Jupiter...\" and \"\u00a9 2020...\" paragraphs (and\n// the blank paragraph that precedes them) that used to follow the\n// \"LOT2002: Biologia Celular (Requisito fraco)\" requirements line.\nconst body = context.document.body;\nbody.paragraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = body.paragraphs.items;\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOT2002: Biologia Celular\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  const toDelete = [];\n  // The three paragraphs right after the requirement line: an empty\n  // paragraph, \"Ver no Jupiter ...\", and \"\u00a9 2020 ...\".\n  for (let i = anchorIndex + 1; i <= anchorIndex + 3 && i < items.length; i++) {\n    toDelete.push(items[i]);\n  }\n  toDelete.forEach((p) => p.delete());\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" and \"\u00a9 2020...\" paragraphs (and\n# the blank paragraph that precedes them) that used to follow the\n# \"LOT2002: Biologia Celular (Requisito fraco)\" requirements line.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*LOT2002: Biologia Celular*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ge 1) {\n    # The three paragraphs right after the requirement line: an empty\n    # paragraph, \"Ver no Jupiter ...\", and \"\u00a9 2020 ...\". Deleting the\n    # paragraph right after the anchor three times removes exactly those,\n    # since each delete shifts the following paragraph into that slot.\n    for ($k = 0; $k -lt 3; $k++) {\n        $p = $d.Paragraphs.Item($anchorIndex + 1)\n        $p.Range.Delete()\n    }\n}\n"}
